# Generate Report for Handoff
# Adds a second tracked file (f8e9b7b9-105a-4e28-ae9c-568a30c60d10) to the
# localization-status workbook: one new row on "Overview", "zh-cn" and
# "de-de", mirroring the existing 67db45b4-... row (same statuses / column
# layout), plus matching hyperlinks.

$wb = $excel.ActiveWorkbook

$newGuid      = "f8e9b7b9-105a-4e28-ae9c-568a30c60d10"
$newGuidMd    = "$newGuid.md"
$newHash      = "e7cbd52349ff8910dbdaba77df21076dc0b9e293"
$mdCommit     = "48ce6e5b62f46b14515eee43e68c64b499b84bbf"
$zhCommit     = "2d905fcff65c4a5da19c8a965bc2a96f08a0027d"
$deCommit     = "65411cd060f565307cebdf09e6e78b09502d9a58"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$newGuidMd"

# ---------------------------------------------------------------------
# Sheet "Overview": add row 3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-30-20 12:30:45"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $newGuidMd) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": add row 3
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhXlf = "$newGuid.$newHash.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"

$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("E3").Value = "2016-03-20 12:30:42"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $newGuidMd) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $mdUrl, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, "", "", $zhXlf) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": add row 3
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deXlf = "$newGuid.$newHash.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("E3").Value = "2016-03-20 12:30:45"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $newGuidMd) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $mdUrl, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, "", "", $deXlf) | Out-Null

Write-Output "Handback report rows added for $newGuid"
